# Apply the edit described by the diff:
#  - Insert a new row for item "محلول ملح" (pushed down to become item #9 at row 15)
#  - Repurpose row 14 (previously item #8 "محلول ملح") to hold a new item
#    "كالونا " (Kalona) with price 15.00 / 15.0000
#  - Update the running total (now row 16) from 289.715 to 304.715
#  - Update the footer timestamp (now row 17) from 2:44 PM to 2:54 PM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original row-14 ("محلول ملح") values before we overwrite them,
# so we can move them down into the freshly inserted row 15.
$origName   = $ws.Range("C14").Value()
$origRatio1 = $ws.Range("H14").Value()
$origLimit  = $ws.Range("N14").Value()
$origPrice  = $ws.Range("P14").Value()

# Insert a new blank row above the current row 15 (the total row),
# shifting the total/footer rows down by one.
$ws.Rows.Item(15).Insert()

# Copy formatting (styles/borders/fills) from row 14 into the new row 15
# so it matches the look of the other item rows.
$ws.Range("A14:Q14").Copy()
$ws.Range("A15:Q15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(15).RowHeight = 24.75

# Merge the new row's cells the same way as the other item rows.
$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

# Helper: write a text value while preserving the cell's numeric display
# format (so numeric-looking text like "15.0000" is stored as text,
# matching the report's existing convention, instead of being coerced
# into a real number).
function Set-TextValue($range, $text) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $fmt
}

# Row 15 (new item #9): move the original "محلول ملح" data down here.
$ws.Range("A15").Value = 9
Set-TextValue $ws.Range("C15") $origName
Set-TextValue $ws.Range("H15") $origRatio1
Set-TextValue $ws.Range("L15") $ws.Range("L14").Value()
Set-TextValue $ws.Range("N15") $origLimit
Set-TextValue $ws.Range("P15") $origPrice
Set-TextValue $ws.Range("Q15") $ws.Range("Q14").Value()

# Row 14 (still item #8): replace with the new item "كالونا ".
Set-TextValue $ws.Range("C14") "كالونا "
Set-TextValue $ws.Range("H14") "0:0"
Set-TextValue $ws.Range("N14") "15.00"
Set-TextValue $ws.Range("P14") "15.0000"

# Update the running total row (now row 16): +15.00 for the new item.
$ws.Range("P16").Value = 304.71499999999997
$ws.Rows.Item(16).RowHeight = 25.5

# Update the footer timestamp (now row 17).
$ws.Range("A17").Value = "Friday, 1 August, 2025 2:54 PM"
